$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D28").Value = "로봇 팔과 카메라 calibration with moveit plugin"
$ws.Range("E28").Value = "https://ropiens.tistory.com/145"

$ws.Range("D37").Value = "[Paper Review] Are Generative Classifiers More Robust to Adversarial Attack?"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1804&mod=document&pageid=1"

$ws.Range("D44").Value = "Forbes IT 기사 리뷰 - Google's Tensor SoC"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/89"

$ws.Range("D46").Value = "[국립암센터] 2021년 08월, 생물정보학(Bioinformatics 채용), 국립암센터 연구소 정규직 연구직 채용공고"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/413"
